# Update the NATMI LR-pairs sheet with the refreshed TPM-based numbers.
# The "Resolving-Mac" sending-cluster rows are retired; in their place the
# "ECs"/"MuSCs" sending-cluster rows gain a new "Resolving-Mac" target row,
# and all the downstream specificity/weight numbers are recomputed. The
# table shrinks from 10 data+header rows (A1:T10) to 9 (A1:T9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data for A2:T9 (the header in row 1 is unchanged).
$data = @(
  @("ECs",   "Amh", "Amhr2", "ECs",           1, 0.3333333333333333, 0.183099,           0.549297,  0.2347163284063858, 0.2347163284063858, 3, 1,                  0.4473893333333334, 1.342168, 0.3318245693495767, 0.3318245693495767, 0.08191653954400001, 0.7372488558960001, 0.07788464459276279, 0.07788464459276277),
  @("ECs",   "Amh", "Amhr2", "FAPs",          1, 0.3333333333333333, 0.183099,           0.549297,  0.2347163284063858, 0.2347163284063858, 2, 0.6666666666666666, 0.416999,           1.250997, 0.3092843375662453, 0.3092843375662453, 0.07635209990099999, 0.687168899109,     0.07259408414715034, 0.07259408414715031),
  @("ECs",   "Amh", "Amhr2", "MuSCs",         1, 0.3333333333333333, 0.183099,           0.549297,  0.2347163284063858, 0.2347163284063858, 3, 1,                  0.4474956666666667, 1.342487, 0.3319034358086358, 0.3319034358086358, 0.081936009071,      0.7374240816390001, 0.07790315583846756, 0.07790315583846756),
  @("ECs",   "Amh", "Amhr2", "Resolving-Mac", 1, 0.3333333333333333, 0.183099,           0.549297,  0.2347163284063858, 0.2347163284063858, 1, 0.3333333333333333, 0.03638666666666666,0.10916,  0.0269876572755421, 0.0269876572755421, 0.00666236228,       0.05996126052,      0.006334443828005128,0.006334443828005126),
  @("MuSCs", "Amh", "Amhr2", "ECs",           3, 1,                  0.5969873333333333,1.790962,  0.7652836715936142, 0.7652836715936142, 3, 1,                  0.4473893333333334, 1.342168, 0.3318245693495767, 0.3318245693495767, 0.2670857650684444,  2.403771885616,     0.2539399247568139, 0.2539399247568139),
  @("MuSCs", "Amh", "Amhr2", "FAPs",          3, 1,                  0.5969873333333333,1.790962,  0.7652836715936142, 0.7652836715936142, 2, 0.6666666666666666, 0.416999,           1.250997, 0.3092843375662453, 0.3092843375662453, 0.2489431210126666,  2.240488089114,     0.236690253419095,  0.236690253419095),
  @("MuSCs", "Amh", "Amhr2", "MuSCs",         3, 1,                  0.5969873333333333,1.790962,  0.7652836715936142, 0.7652836715936142, 3, 1,                  0.4474956666666667, 1.342487, 0.3319034358086358, 0.3319034358086358, 0.2671492447215555,  2.404343202494,     0.2540002799701683, 0.2540002799701683),
  @("MuSCs", "Amh", "Amhr2", "Resolving-Mac", 3, 1,                  0.5969873333333333,1.790962,  0.7652836715936142, 0.7652836715936142, 1, 0.3333333333333333, 0.03638666666666666,0.10916,  0.0269876572755421, 0.0269876572755421, 0.02172237910222222, 0.19550141192,      0.02065321344753698,0.02065321344753697)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $data[$i]
  $r = $startRow + $i
  for ($c = 0; $c -lt $row.Count; $c++) {
    $ws.Cells.Item($r, $c + 1).Value = $row[$c]
  }
}

# The old row 10 (former "Resolving-Mac" sending-cluster data) no longer
# exists; the used range shrinks to A1:T9.
$ws.Range("A10:T10").Clear()
